$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 105, shifting existing rows
# 105-189 down to 106-190 (formatting is inherited from the row above,
# which keeps the date style on column D).
$ws.Rows(105).Insert()

# Populate the newly inserted row 105 with the new record's data.
$ws.Cells.Item(105, 1).Value2  = 4
$ws.Cells.Item(105, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(105, 3).Value2  = "Los Lagos"
$ws.Cells.Item(105, 4).Value2  = 44574
$ws.Cells.Item(105, 5).Value2  = 10
$ws.Cells.Item(105, 6).Value2  = "Fruta"
$ws.Cells.Item(105, 7).Value2  = 100108
$ws.Cells.Item(105, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(105, 9).Value2  = 100108005
$ws.Cells.Item(105, 10).Value2 = "Piña"
$ws.Cells.Item(105, 11).Value2 = "Caramelo"
$ws.Cells.Item(105, 12).Value2 = "Tercera"
$ws.Cells.Item(105, 13).Value2 = 120
$ws.Cells.Item(105, 14).Value2 = 19000
$ws.Cells.Item(105, 15).Value2 = 20000
$ws.Cells.Item(105, 16).Value2 = 19500
$ws.Cells.Item(105, 17).Value2 = "`$/caja 16 unidades"
$ws.Cells.Item(105, 18).Value2 = "Ecuador"
$ws.Cells.Item(105, 19).Value2 = 1219
$ws.Cells.Item(105, 20).Value2 = 16
